$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.955.18'
$ws.Range("D3").Value = '3.146.80'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.22%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.141.98'
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.39'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.52%  '
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").Value = '3.662.74'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("D17").Value = '63.972.41'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").Value = '3.132.10'
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '488.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.24%  '
$ws.Range("E25").Value = '  -3.59%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.68%  '
$ws.Range("E31").Value = '  +3.76%  '
$ws.Range("E32").Value = '  -5.29%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("D38").Value = '0.0₃0748'
$ws.Range("E38").Value = '  -4.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '434.55'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0397'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("D44").Value = '2.933.05'
$ws.Range("E44").Value = '  +2.91%  '
$ws.Range("E45").Value = '  -2.73%  '
$ws.Range("E46").Value = '  -5.11%  '
$ws.Range("E47").Value = '  -1.12%  '
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.20%  '
